$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 37.22222
$ws.Range("I9").Value = 24.857143
$ws.Range("J9").Value = 80.5
$ws.Range("K9").Value = 24.857143
$ws.Range("L9").Value = 80.5
$ws.Range("M9").Value = 144.142857
$ws.Range("N9").Value = -418.5
$ws.Range("H15").Value = 1608.2759
$ws.Range("I15").Value = 1608.2759
$ws.Range("K15").Value = 4824.8277
$ws.Range("M15").Value = -4655.8277
$ws.Range("H17").Value = 3775.8965
$ws.Range("J17").Value = 3900
$ws.Range("L17").Value = 11700
$ws.Range("N17").Value = -12036
$ws.Range("H55").Value = 393.25
$ws.Range("J55").Value = 479
$ws.Range("L55").Value = 479
$ws.Range("N55").Value = -907
$ws.Range("H69").Value = 3002
$ws.Range("I69").Value = 2671
$ws.Range("J69").Value = 3995
$ws.Range("K69").Value = 8013
$ws.Range("L69").Value = 11985
$ws.Range("M69").Value = -7139
$ws.Range("N69").Value = -13733
$ws.Range("H72").Value = 3002
$ws.Range("I72").Value = 2671
$ws.Range("J72").Value = 3995
$ws.Range("K72").Value = 24039
$ws.Range("L72").Value = 35955
$ws.Range("M72").Value = -19671
$ws.Range("N72").Value = -44691
$ws.Range("H113").Value = 4212
$ws.Range("I113").Value = 3005
$ws.Range("J113").Value = 4453.4
$ws.Range("K113").Value = 3005
$ws.Range("L113").Value = 4453.4
$ws.Range("M113").Value = 249
$ws.Range("N113").Value = -10961.4

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2859
$ws.Range("I2").Value = 3450
$ws.Range("J2").Value = 1086
$ws.Range("K2").Value = 3450
$ws.Range("L2").Value = 1086
$ws.Range("M2").Value = -3337
$ws.Range("N2").Value = -1312
$ws.Range("H110").Value = 3181.3333
$ws.Range("I110").Value = 635.2
$ws.Range("K110").Value = 635.2
$ws.Range("M110").Value = 1409.8
$ws.Range("H116").Value = 2859
$ws.Range("I116").Value = 3450
$ws.Range("J116").Value = 1086
$ws.Range("K116").Value = 3450
$ws.Range("L116").Value = 1086
$ws.Range("M116").Value = -1156
$ws.Range("N116").Value = -5674
$ws.Range("H122").Value = 2045.4445
$ws.Range("I122").Value = 2045.4445
$ws.Range("K122").Value = 6136.333500000001
$ws.Range("M122").Value = -3686.333500000001
$ws.Range("H132").Value = 2135.625
$ws.Range("I132").Value = 2135.625
$ws.Range("K132").Value = 6406.875
$ws.Range("M132").Value = -3876.875

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2859
$ws.Range("I3").Value = 3450
$ws.Range("J3").Value = 1086
$ws.Range("K3").Value = 3450
$ws.Range("L3").Value = 1086
$ws.Range("M3").Value = -3336
$ws.Range("N3").Value = -1314

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("H58").Value = 7492.7856
$ws.Range("I58").Value = 4481.727
$ws.Range("K58").Value = 4481.727
$ws.Range("M58").Value = -4278.727
$ws.Range("H86").Value = 4597.7
$ws.Range("I86").Value = 4424.143
$ws.Range("K86").Value = 4424.143
$ws.Range("M86").Value = -3301.143
$ws.Range("H89").Value = 4597.7
$ws.Range("I89").Value = 4424.143
$ws.Range("K89").Value = 22120.715
$ws.Range("M89").Value = -16504.715
$ws.Range("H134").Value = 3054.111
$ws.Range("I134").Value = 2550.8462
$ws.Range("J134").Value = 4362.6
$ws.Range("K134").Value = 7652.5386
$ws.Range("L134").Value = 13087.8
$ws.Range("M134").Value = -5117.5386
$ws.Range("N134").Value = -18157.8
$ws.Range("H136").Value = 7492.7856
$ws.Range("I136").Value = 4481.727
$ws.Range("K136").Value = 13445.181
$ws.Range("M136").Value = -10895.181
$ws.Range("M57").ClearContents()

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 215.66667
$ws.Range("J12").Value = 227.09091
$ws.Range("L12").Value = 681.27273
$ws.Range("N12").Value = -1027.27273
$ws.Range("H32").Value = 787.5
$ws.Range("I32").Value = 750
$ws.Range("K32").Value = 2250
$ws.Range("M32").Value = -1967
$ws.Range("H63").Value = 2000
$ws.Range("I63").Value = 2000
$ws.Range("K63").Value = 6000
$ws.Range("M63").Value = -5251
$ws.Range("H66").Value = 2000
$ws.Range("I66").Value = 2000
$ws.Range("K66").Value = 18000
$ws.Range("M66").Value = -14256
$ws.Range("H114").Value = 0
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").Value = 0
$ws.Range("H122").Value = 570.2222
$ws.Range("I122").Value = 565.3333
$ws.Range("K122").Value = 5087.9997
$ws.Range("M122").Value = -2637.9997
$ws.Range("M114").ClearContents()

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("H49").Value = 24500
$ws.Range("I49").Value = 24000
$ws.Range("K49").Value = 24000
$ws.Range("M49").Value = -23816
$ws.Range("H122").Value = 10451320
$ws.Range("I122").Value = 12540384
$ws.Range("K122").Value = 37621152
$ws.Range("M122").Value = -37618702
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("M126").ClearContents()

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2143.6
$ws.Range("I7").Value = 1554.5
$ws.Range("K7").Value = 1554.5
$ws.Range("M7").Value = -1442.5
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("H126").Value = 2143.6
$ws.Range("I126").Value = 1554.5
$ws.Range("K126").Value = 4663.5
$ws.Range("M126").Value = -2193.5
$ws.Range("N98").ClearContents()

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = 0
$ws.Range("M132").ClearContents()
